$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197, shifting existing rows 197-269 down to 198-270
$ws.Rows.Item(197).Insert()

# Populate the new row 197 with the new data record
$ws.Range("A197").Value = 1
$ws.Range("B197").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C197").Value = "Arica y Parinacota"
$ws.Range("D197").Value = 45146
$ws.Range("E197").Value = 15
$ws.Range("F197").Value = "Fruta"
$ws.Range("G197").Value = 100104
$ws.Range("H197").Value = "Frutos de pepita"
$ws.Range("I197").Value = 100104002
$ws.Range("J197").Value = "Manzana"
$ws.Range("K197").Value = "Fuji royal"
$ws.Range("L197").Value = "Calibre 100"
$ws.Range("M197").Value = 300
$ws.Range("N197").Value = 20000
$ws.Range("O197").Value = 22000
$ws.Range("P197").Value = 21000
$ws.Range("Q197").Value = "`$/caja 18 kilos embalada"
$ws.Range("R197").Value = "Región de O'Higgins"
$ws.Range("S197").Value = 1167
$ws.Range("T197").Value = 18
